$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row 9 with the new log entry
$ws.Range("A9").Value = "2/13/2020jaclemon"
$ws.Range("C9").Value = "15 minutes"
$ws.Range("D9").Value = "Used Clion to recognize illegal command line args"

# Move selection to D10, matching the post-edit cursor position
$ws.Range("D10").Select()
